# Rename the header cell A1 from "Name" to "Tag", and touch the formatting
# on the A1:A2 column so both cells pick up an explicit (new) cell style
# instead of sharing the sheet's implicit default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 held the shared string "Name" -> rename it to "Tag".
$ws.Range("A1").Value = "Tag"

# Re-apply the "Normal" style explicitly across A1:A2 so both cells get a
# fresh, explicit cellXf entry (style index 1) instead of implicitly
# referencing the default style (index 0).
$ws.Range("A1:A2").Style = "Normal"

# Leave the selection on A2, matching where the cursor ended up.
[void]$ws.Range("A2").Select()
